# Update res_bus/vm_pu.xlsx values for the "case with 380 kV" run.
# Columns B-F and I-M (bus voltage magnitudes, pu) are refreshed for rows 2-25;
# columns A, G, H, N are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 1.02
$ws.Range("C2").Value2 = 1.040555794334618
$ws.Range("D2").Value2 = 1.047096738552105
$ws.Range("E2").Value2 = 1.048319230427386
$ws.Range("F2").Value2 = 1.057694054655205
$ws.Range("I2").Value2 = 1.035998129261776
$ws.Range("J2").Value2 = 1.045641764919756
$ws.Range("K2").Value2 = 1.049860330828642
$ws.Range("L2").Value2 = 1.051079407225533
$ws.Range("M2").Value2 = 1.060428323745304

$ws.Range("B3").Value2 = 1.02
$ws.Range("C3").Value2 = 1.041556995282184
$ws.Range("D3").Value2 = 1.047870163236246
$ws.Range("E3").Value2 = 1.049196448301861
$ws.Range("F3").Value2 = 1.058624907041211
$ws.Range("I3").Value2 = 1.036153347773119
$ws.Range("J3").Value2 = 1.046288373796804
$ws.Range("K3").Value2 = 1.050445630435419
$ws.Range("L3").Value2 = 1.051768477786647
$ws.Range("M3").Value2 = 1.061172765775723

$ws.Range("B4").Value2 = 1.02
$ws.Range("C4").Value2 = 1.042205390041193
$ws.Range("D4").Value2 = 1.048370878913516
$ws.Range("E4").Value2 = 1.049764887251109
$ws.Range("F4").Value2 = 1.059228030536945
$ws.Range("I4").Value2 = 1.036252507435273
$ws.Range("J4").Value2 = 1.046706721272317
$ws.Range("K4").Value2 = 1.05082395628745
$ws.Range("L4").Value2 = 1.052214533204774
$ws.Range("M4").Value2 = 1.061654639138753

$ws.Range("B5").Value2 = 1.02
$ws.Range("C5").Value2 = 1.042478106270181
$ws.Range("D5").Value2 = 1.048581440270978
$ws.Range("E5").Value2 = 1.050004054122145
$ws.Range("F5").Value2 = 1.059481773805033
$ws.Range("I5").Value2 = 1.036293887905323
$ws.Range("J5").Value2 = 1.046882581263618
$ws.Range("K5").Value2 = 1.050982907191634
$ws.Range("L5").Value2 = 1.052402097155591
$ws.Range("M5").Value2 = 1.061857258181567

$ws.Range("B6").Value2 = 1.02
$ws.Range("C6").Value2 = 1.042523904188341
$ws.Range("D6").Value2 = 1.048616797957882
$ws.Range("E6").Value2 = 1.050044222698744
$ws.Range("F6").Value2 = 1.059524389528642
$ws.Range("I6").Value2 = 1.036300817894674
$ws.Range("J6").Value2 = 1.046912108124407
$ws.Range("K6").Value2 = 1.051009590013652
$ws.Range("L6").Value2 = 1.052433592416828
$ws.Range("M6").Value2 = 1.061891281093097

$ws.Range("B7").Value2 = 1.02
$ws.Range("C7").Value2 = 1.042209033574746
$ws.Range("D7").Value2 = 1.048373692208673
$ws.Range("E7").Value2 = 1.049768082245361
$ws.Range("F7").Value2 = 1.059231420321397
$ws.Range("I7").Value2 = 1.036253061567139
$ws.Range("J7").Value2 = 1.046709071175477
$ws.Range("K7").Value2 = 1.050826080580743
$ws.Range("L7").Value2 = 1.052217039278698
$ws.Range("M7").Value2 = 1.06165734639059

$ws.Range("B8").Value2 = 1.02
$ws.Range("C8").Value2 = 1.040894040359946
$ws.Range("D8").Value2 = 1.047358066789423
$ws.Range("E8").Value2 = 1.048615519708758
$ws.Range("F8").Value2 = 1.058008474192383
$ws.Range("I8").Value2 = 1.036050850206405
$ws.Range("J8").Value2 = 1.04586029957624
$ws.Range("K8").Value2 = 1.050058218225935
$ws.Range("L8").Value2 = 1.051312243831236
$ws.Range("M8").Value2 = 1.060679875627089

$ws.Range("B9").Value2 = 1.02
$ws.Range("C9").Value2 = 1.038581107598058
$ws.Range("D9").Value2 = 1.045570440869252
$ws.Range("E9").Value2 = 1.046590895187545
$ws.Range("F9").Value2 = 1.05585966999995
$ws.Range("I9").Value2 = 1.035684769595194
$ws.Range("J9").Value2 = 1.044364299665574
$ws.Range("K9").Value2 = 1.048702111481478
$ws.Range("L9").Value2 = 1.04971930680116
$ws.Range("M9").Value2 = 1.058958795487696

$ws.Range("B10").Value2 = 1.02
$ws.Range("C10").Value2 = 1.037042057483732
$ws.Range("D10").Value2 = 1.044380136999011
$ws.Range("E10").Value2 = 1.045245480319907
$ws.Range("F10").Value2 = 1.054431368238821
$ws.Range("I10").Value2 = 1.035434185101928
$ws.Range("J10").Value2 = 1.043366780280379
$ws.Range("K10").Value2 = 1.04779605912055
$ws.Range("L10").Value2 = 1.048658367880479
$ws.Range("M10").Value2 = 1.057812380908639

$ws.Range("B11").Value2 = 1.02
$ws.Range("C11").Value2 = 1.036376329879835
$ws.Range("D11").Value2 = 1.0438650816192
$ws.Range("E11").Value2 = 1.044663944218564
$ws.Range("F11").Value2 = 1.053813918111904
$ws.Range("I11").Value2 = 1.035324136490078
$ws.Range("J11").Value2 = 1.042934810156772
$ws.Range("K11").Value2 = 1.047403270820958
$ws.Range("L11").Value2 = 1.048199223480447
$ws.Range("M11").Value2 = 1.057316214899727

$ws.Range("B12").Value2 = 1.02
$ws.Range("C12").Value2 = 1.03612915352235
$ws.Range("D12").Value2 = 1.043673821584395
$ws.Range("E12").Value2 = 1.044448092875635
$ws.Range("F12").Value2 = 1.053584723232796
$ws.Range("I12").Value2 = 1.035283028058852
$ws.Range("J12").Value2 = 1.042774352404566
$ws.Range("K12").Value2 = 1.047257303352428
$ws.Range("L12").Value2 = 1.048028715184055
$ws.Range("M12").Value2 = 1.057131953726951

$ws.Range("B13").Value2 = 1.02
$ws.Range("C13").Value2 = 1.036182168949806
$ws.Range("D13").Value2 = 1.043714845025951
$ws.Range("E13").Value2 = 1.044494386594967
$ws.Range("F13").Value2 = 1.053633879341972
$ws.Range("I13").Value2 = 1.035291856424111
$ws.Range("J13").Value2 = 1.042808771341458
$ws.Range("K13").Value2 = 1.047288616950105
$ws.Range("L13").Value2 = 1.048065288040608
$ws.Range("M13").Value2 = 1.057171476683828

$ws.Range("B14").Value2 = 1.02
$ws.Range("C14").Value2 = 1.036355896069063
$ws.Range("D14").Value2 = 1.043849270886601
$ws.Range("E14").Value2 = 1.044646098668354
$ws.Range("F14").Value2 = 1.053794969651255
$ws.Range("I14").Value2 = 1.035320743173649
$ws.Range("J14").Value2 = 1.04292154676227
$ws.Range("K14").Value2 = 1.047391206494409
$ws.Range("L14").Value2 = 1.048185128423127
$ws.Range("M14").Value2 = 1.057300983053313

$ws.Range("B15").Value2 = 1.02
$ws.Range("C15").Value2 = 1.036462948950305
$ws.Range("D15").Value2 = 1.043932102351352
$ws.Range("E15").Value2 = 1.044739594331465
$ws.Range("F15").Value2 = 1.053894243097149
$ws.Range("I15").Value2 = 1.03533851059465
$ws.Range("J15").Value2 = 1.042991030808274
$ws.Range("K15").Value2 = 1.047454406262691
$ws.Range("L15").Value2 = 1.048258971151658
$ws.Range("M15").Value2 = 1.057380781132963

$ws.Range("B16").Value2 = 1.02
$ws.Range("C16").Value2 = 1.037086254252651
$ws.Range("D16").Value2 = 1.044414327096469
$ws.Range("E16").Value2 = 1.0452840969151
$ws.Range("F16").Value2 = 1.054472367840209
$ws.Range("I16").Value2 = 1.035441456178543
$ws.Range("J16").Value2 = 1.043395447993408
$ws.Range("K16").Value2 = 1.047822117541241
$ws.Range("L16").Value2 = 1.048688845112142
$ws.Range("M16").Value2 = 1.057845314939697

$ws.Range("B17").Value2 = 1.02
$ws.Range("C17").Value2 = 1.037477422735616
$ws.Range("D17").Value2 = 1.044716909627804
$ws.Range("E17").Value2 = 1.045625927591253
$ws.Range("F17").Value2 = 1.054835282568614
$ws.Range("I17").Value2 = 1.035505618218172
$ws.Range("J17").Value2 = 1.043649118723252
$ws.Range("K17").Value2 = 1.048052650187129
$ws.Range("L17").Value2 = 1.048958561062505
$ws.Range("M17").Value2 = 1.058136769475463

$ws.Range("B18").Value2 = 1.02
$ws.Range("C18").Value2 = 1.037705651344661
$ws.Range("D18").Value2 = 1.044893434974385
$ws.Range("E18").Value2 = 1.045825411656417
$ws.Range("F18").Value2 = 1.055047062338043
$ws.Range("I18").Value2 = 1.035542893857263
$ws.Range("J18").Value2 = 1.043797076879283
$ws.Range("K18").Value2 = 1.048187071372888
$ws.Range("L18").Value2 = 1.049115905789424
$ws.Range("M18").Value2 = 1.058306792970869

$ws.Range("B19").Value2 = 1.02
$ws.Range("C19").Value2 = 1.037783482687455
$ws.Range("D19").Value2 = 1.044953631296755
$ws.Range("E19").Value2 = 1.045893447495361
$ws.Range("F19").Value2 = 1.055119290308234
$ws.Range("I19").Value2 = 1.035555578600735
$ws.Range("J19").Value2 = 1.043847526135637
$ws.Range("K19").Value2 = 1.048232897928522
$ws.Range("L19").Value2 = 1.049169560313229
$ws.Range("M19").Value2 = 1.058364770457326

$ws.Range("B20").Value2 = 1.02
$ws.Range("C20").Value2 = 1.037435447160507
$ws.Range("D20").Value2 = 1.044684441852003
$ws.Range("E20").Value2 = 1.045589242041043
$ws.Range("F20").Value2 = 1.054796335133186
$ws.Range("I20").Value2 = 1.035498749646175
$ws.Range("J20").Value2 = 1.043621902630112
$ws.Range("K20").Value2 = 1.048027920822896
$ws.Range("L20").Value2 = 1.048929620606764
$ws.Range("M20").Value2 = 1.058105496775853

$ws.Range("B21").Value2 = 1.02
$ws.Range("C21").Value2 = 1.036304734871964
$ws.Range("D21").Value2 = 1.043809684307356
$ws.Range("E21").Value2 = 1.044601418898869
$ws.Range("F21").Value2 = 1.05374752832033
$ws.Range("I21").Value2 = 1.035312243131987
$ws.Range("J21").Value2 = 1.042888337337843
$ws.Range("K21").Value2 = 1.047360998311253
$ws.Range("L21").Value2 = 1.048149837350427
$ws.Range("M21").Value2 = 1.057262845626181

$ws.Range("B22").Value2 = 1.02
$ws.Range("C22").Value2 = 1.035594416209079
$ws.Range("D22").Value2 = 1.043260005274548
$ws.Range("E22").Value2 = 1.043981244565766
$ws.Range("F22").Value2 = 1.053088991096057
$ws.Range("I22").Value2 = 1.035193640065215
$ws.Range("J22").Value2 = 1.042427088284076
$ws.Range("K22").Value2 = 1.046941282359052
$ws.Range("L22").Value2 = 1.047659779174776
$ws.Range("M22").Value2 = 1.056733251852371

$ws.Range("B23").Value2 = 1.02
$ws.Range("C23").Value2 = 1.035970911998998
$ws.Range("D23").Value2 = 1.043551370174906
$ws.Range("E23").Value2 = 1.044309924124005
$ws.Range("F23").Value2 = 1.053438009429338
$ws.Range("I23").Value2 = 1.035256640560336
$ws.Range("J23").Value2 = 1.042671607454195
$ws.Range("K23").Value2 = 1.04716381887817
$ws.Range("L23").Value2 = 1.047919546879043
$ws.Range("M23").Value2 = 1.057013978851475

$ws.Range("B24").Value2 = 1.02
$ws.Range("C24").Value2 = 1.037454413910139
$ws.Range("D24").Value2 = 1.044699112537156
$ws.Range("E24").Value2 = 1.045605818353146
$ws.Range("F24").Value2 = 1.054813933500863
$ws.Range("I24").Value2 = 1.035501853718607
$ws.Range("J24").Value2 = 1.043634200422552
$ws.Range("K24").Value2 = 1.048039095096028
$ws.Range("L24").Value2 = 1.048942697478975
$ws.Range("M24").Value2 = 1.058119627491718

$ws.Range("B25").Value2 = 1.02
$ws.Range("C25").Value2 = 1.039178547463047
$ws.Range("D25").Value2 = 1.0460323357271
$ws.Range("E25").Value2 = 1.047113550503223
$ws.Range("F25").Value2 = 1.056414446708348
$ws.Range("I25").Value2 = 1.035780563514513
$ws.Range("J25").Value2 = 1.044751088203617
$ws.Range("K25").Value2 = 1.0490530511375
$ws.Range("L25").Value2 = 1.050130943898067
$ws.Range("M25").Value2 = 1.059403569418627
